# "unidades de almacenamiento y configuracion lan"
# Fills in the (previously empty) "unidades" sheet with a storage-units /
# bandwidth-units reference table, then updates the view/selection state
# on both sheets so "unidades" becomes the active tab.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("modelos")
$ws2 = $wb.Worksheets.Item("unidades")

# ---------------------------------------------------------------------
# 1) Set the (soon to be former) selection on "modelos" before we move
#    away from it, so its sheetView keeps a C2:C4 selection and loses
#    the "active tab" flag once "unidades" is activated below.
# ---------------------------------------------------------------------
$ws1.Range("C2:C4").Select()

# ---------------------------------------------------------------------
# 2) Populate "unidades" - write cells in the same order the original
#    author did, so that the shared-string table comes out in the same
#    sequence.
# ---------------------------------------------------------------------
$ws2.Range("A1").Value = "bit"
$ws2.Range("B1").Value = "Binary Digit"
$ws2.Range("C1").Value = "U. Minima"

$ws2.Range("A3").Value = "ALMACENAMIENTO"
$ws2.Range("E3").Value = "ANCHO DE BANDA"

$ws2.Range("A4").Value = "Unidad"
$ws2.Range("B4").Value = "Equivalencia"
$ws2.Range("C4").Value = "Exponente"

$ws2.Range("A5").Value = "Byte"
$ws2.Range("B5").Value = "8 bits"
$ws2.Range("C5").Value = "10^0"

$ws2.Range("A6").Value = "KiloByte"

$ws2.Range("C6").Value  = "10^3"
$ws2.Range("C7").Value  = "10^6"
$ws2.Range("C8").Value  = "10^9"
$ws2.Range("C9").Value  = "10^12"
$ws2.Range("C10").Value = "10^15"
$ws2.Range("C11").Value = "10^18"
$ws2.Range("C12").Value = "10^21"
$ws2.Range("C13").Value = "10^24"
$ws2.Range("C14").Value = "10^27"
$ws2.Range("C15").Value = "10^30"
$ws2.Range("C16").Value = "10^33"

$ws2.Range("A7").Value = "MegaByte"
$ws2.Range("B6").Value = "1000B"
$ws2.Range("B7").Value = "1000KB"

$ws2.Range("A8").Value = "GigaByte"
$ws2.Range("B8").Value = "1000MB"

$ws2.Range("A9").Value = "TeraByte"
$ws2.Range("B9").Value = "1000GB"

$ws2.Range("A10").Value = "PetaByte"
$ws2.Range("A11").Value = "ExaByte"
$ws2.Range("A12").Value = "ZettaByte"
$ws2.Range("A13").Value = "YottaByte"
$ws2.Range("A14").Value = "BrontoByte"
$ws2.Range("A15").Value = "GeopByte"
$ws2.Range("A16").Value = "SaganByte"

$ws2.Range("B10").Value = "1000TB"
$ws2.Range("B11").Value = "1000PB"
$ws2.Range("B12").Value = "1000XB"
$ws2.Range("B13").Value = "1000ZB"
$ws2.Range("B14").Value = "1000YB"
$ws2.Range("B15").Value = "1000BB"
$ws2.Range("B16").Value = "1000GeB"

$ws2.Range("D6").Value = "KiB"
$ws2.Range("D8").Value = "GiB"
$ws2.Range("D7").Value = "MiB"
$ws2.Range("D9").Value = "TiB"

# Mirrored "ANCHO DE BANDA" header row (reuses the same 3 strings).
$ws2.Range("E4").Value = "Unidad"
$ws2.Range("F4").Value = "Equivalencia"
$ws2.Range("G4").Value = "Exponente"

# ---------------------------------------------------------------------
# 3) Formatting: merge + center the two section titles, and left-indent
#    the data/header cells.
# ---------------------------------------------------------------------
$ws2.Range("A3:C3").Merge()
$ws2.Range("A3:C3").HorizontalAlignment = -4108   # xlCenter
$ws2.Range("E3:G3").Merge()
$ws2.Range("E3:G3").HorizontalAlignment = -4108   # xlCenter

$ws2.Range("D4:G4").HorizontalAlignment = -4108   # xlCenter

$left = $ws2.Range("A4:C16")
$left.HorizontalAlignment = -4131   # xlLeft
$left.IndentLevel = 1

$colD = $ws2.Range("D6:D9")
$colD.HorizontalAlignment = -4131   # xlLeft
$colD.IndentLevel = 1

# Column B is the widest column on this sheet (bestFit-style sizing).
$ws2.Columns.Item(2).ColumnWidth = 12.59

# ---------------------------------------------------------------------
# 4) Add the thin spacer column on "modelos" (col E) that shows up in
#    the diff's <cols> list.
# ---------------------------------------------------------------------
$ws1.Columns.Item(5).ColumnWidth = 1.74

# ---------------------------------------------------------------------
# 5) View state: make "unidades" the active sheet/tab, restore its zoom,
#    scroll position and selection.
# ---------------------------------------------------------------------
$ws2.Activate()
$excel.ActiveWindow.Zoom = 145
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$ws2.Range("D6:D9").Select()

Write-Host "done"
